$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '58.300.45'
$ws.Range('E2').Value = '  -0.77%  '

# Row 3
$ws.Range('D3').Value = '2.483.43'
$ws.Range('E3').Value = '  -1.41%  '

# Row 4
$ws.Range('E4').Value = '  +0.09%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '521.51'
$ws.Range('E5').Value = '  -2.66%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.20'
$ws.Range('E6').Value = '  -3.18%  '

# Row 7
$ws.Range('E7').Value = '  +0.16%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.560'
$ws.Range('E8').Value = '  -1.09%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0998'
$ws.Range('E9').Value = '  -0.98%  '

# Row 10
$ws.Range('E10').Value = '  -0.26%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.38'
$ws.Range('E11').Value = '  +0.86%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.345'
$ws.Range('E12').Value = '  -0.26%  '

# Row 13
$ws.Range('D13').Value = '2.923.05'
$ws.Range('E13').Value = '  -0.42%  '

# Row 14
$ws.Range('D14').Value = '58.216.58'
$ws.Range('E14').Value = '  -0.78%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.48'
$ws.Range('E15').Value = '  -1.69%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000137'
$ws.Range('E16').Value = '  -1.24%  '

# Row 17
$ws.Range('D17').Value = '2.484.06'
$ws.Range('E17').Value = '  -0.96%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.89'
$ws.Range('E18').Value = '  -1.50%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.20'
$ws.Range('E19').Value = '  -1.56%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '320.19'
$ws.Range('E20').Value = '  -0.83%  '

# Row 21
$ws.Range('E21').Value = '  +0.05%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.78'
$ws.Range('E22').Value = '  -2.23%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.34'
$ws.Range('E23').Value = '  -1.54%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.412'
$ws.Range('E24').Value = '  -1.77%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  +0.10%  '

# Row 26
$ws.Range('E26').Value = '  -2.95%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.39'
$ws.Range('E27').Value = '  -2.06%  '

# Row 28
$ws.Range('D28').Value = '0.0₃0760'
$ws.Range('E28').Value = '  -0.72%  '

# Row 29
$ws.Range('E29').Value = '  -2.87%  '

# Row 30
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.36'
$ws.Range('E30').Value = '  -4.82%  '

# Row 31
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '166.56'
$ws.Range('E31').Value = '  -0.02%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.18'
$ws.Range('E32').Value = '  +1.63%  '

# Row 33
$ws.Range('E33').Value = '  +0.07%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  +0.01%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.16'
$ws.Range('E35').Value = '  -1.46%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.33'
$ws.Range('E36').Value = '  -8.68%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.00'
$ws.Range('E37').Value = '  -1.92%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.49'
$ws.Range('E38').Value = '  -3.06%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.796'
$ws.Range('E39').Value = '  -1.95%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.50'
$ws.Range('E40').Value = '  -2.80%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '277.49'
$ws.Range('E41').Value = '  -2.36%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.08'
$ws.Range('E42').Value = '  -2.66%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.598'
$ws.Range('E43').Value = '  -0.94%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '127.66'
$ws.Range('E44').Value = '  -2.47%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0909'
$ws.Range('E45').Value = '  -1.50%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0491'
$ws.Range('E46').Value = '  -2.75%  '

# Row 47
$ws.Range('E47').Value = '  -2.24%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '17.22'
$ws.Range('E48').Value = '  +0.22%  '

# Row 49
$ws.Range('D49').Value = '1.745.95'
$ws.Range('E49').Value = '  -1.17%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.974'
$ws.Range('E50').Value = '  -0.88%  '

# Row 51
$ws.Range('E51').Value = '  -0.98%  '
